$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.00006268748566506943
$ws.Range("C2").Value = 0.001599156226475884
$ws.Range("D2").Value = 0.0000430553141705925
$ws.Range("B3").Value = 0.00000009107660048357502
$ws.Range("C3").Value = 0.00000176407894514341
$ws.Range("D3").Value = 0.0000001275072492035179
$ws.Range("B4").Value = 0.000000002296691947734075
$ws.Range("C4").Value = 0.000008280819201900158
$ws.Range("D4").Value = 0.00000002038774660206855
$ws.Range("B5").Value = 0.00003131014289015077
$ws.Range("C5").Value = 0.000424092825120681
$ws.Range("D5").Value = 0.00001038567922523326
$ws.Range("B6").Value = 0.0000000416423208693395
$ws.Range("C6").Value = 0.000001062294590070678
$ws.Range("D6").Value = 0.00000002860097503276648
$ws.Range("B7").Value = 0.000004172967717153142
$ws.Range("C7").Value = 0.01007666458463063
$ws.Range("D7").Value = 0.00003704345667188136
$ws.Range("B8").Value = 0.003141038826697695
$ws.Range("C8").Value = 0.1255649392405758
$ws.Range("D8").Value = 0.001554776877810582
$ws.Range("B9").Value = 0.00004552381112477022
$ws.Range("C9").Value = 0.001404903906873756
$ws.Range("D9").Value = 0.002599266092119024
$ws.Range("B10").Value = 0.008508961102393187
$ws.Range("C10").Value = 0.07507947120211611
$ws.Range("D10").Value = 0.0002123632777966122
$ws.Range("B11").Value = 0.000002413879530038798
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0.00004223355209376223
